$d = $word.ActiveDocument
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Tipovi</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p/><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="A31515"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="A31515"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>api</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="A31515"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>/portal/types</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="A31515"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="A31515"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>getalltypes</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="A31515"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="A31515"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>inserttypes</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="A31515"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="A31515"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>updatetypes</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="A31515"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Gradovi</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p/><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="A31515"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="A31515"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>api</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="A31515"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>/portal/cities</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="A31515"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="A31515"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>getallcities</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="A31515"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="A31515"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>insertcities</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="A31515"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="A31515"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>updatecities</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="A31515"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Osobe</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="C00000"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="C00000"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>api</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="C00000"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>/portal/persons</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="C00000"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="C00000"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>getallpersons</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="C00000"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="C00000"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>insertperson</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="C00000"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="C00000"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>updateperson</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p/><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Apartmani</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="C00000"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="C00000"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>api</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="C00000"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>/portal/apartments</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="C00000"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="C00000"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>getallapartments</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="C00000"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="C00000"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>insertapartment</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="C00000"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/><w:color w:val="C00000"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>updateapartment</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$d.Content.InsertXML($xml)
